$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks from D2, D3, F3 but keep the Hyperlink cell style (s="1")
# that they already carry.
$ws.Hyperlinks.Delete()

# Update row 2 values: player name and photo link change.
$ws.Range("B2").Value = "Tom"
$ws.Range("D2").Value = "https://www.flickr.com/photos/jakubz/30287928287/in/feed"

# Clear row 2 trailing cells (E2, F2) that are no longer populated.
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

# Row 3 keeps only D3 and F3 as empty (but still styled) cells; everything
# else in row 3 is cleared.
$ws.Range("A3:H3").ClearContents()

# Rows 4 and 5 are removed entirely from the sheet.
$ws.Range("A4:H5").Delete(-4162)

# Restore the hyperlink-style formatting (font/underline) on D3 and F3 even
# though they no longer carry an actual hyperlink.
$ws.Range("D3").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"

$ws.Range("D2").Select()
